# The commit inserts one new weekly price-report row for Cilantro
# (Femacal de La Calera) at sheet row 156, pushing the existing rows
# 156-245 down to 157-246 (dimension grows from A1:R245 to A1:R246).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 156, shifting everything
# below it (including all formatting) down by one row.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new record.
$ws.Cells.Item(156, 1).Value  = 3
$ws.Cells.Item(156, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(156, 3).Value  = "Coquimbo"
$ws.Cells.Item(156, 4).Value  = 44518
$ws.Cells.Item(156, 5).Value  = 5
$ws.Cells.Item(156, 6).Value  = 100112040
$ws.Cells.Item(156, 7).Value  = "Cilantro"
$ws.Cells.Item(156, 8).Value  = "Sin especificar"
$ws.Cells.Item(156, 9).Value  = "Primera"
$ws.Cells.Item(156, 10).Value = 280
$ws.Cells.Item(156, 11).Value = 2000
$ws.Cells.Item(156, 12).Value = 2500
$ws.Cells.Item(156, 13).Value = 2214
$ws.Cells.Item(156, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(156, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(156, 16).Value = 738
$ws.Cells.Item(156, 17).Value = 3
$ws.Cells.Item(156, 18).Value = "Hortaliza"
